# Insert a new data row above row 135 (pushes the existing rows 135-202
# down to 136-203, extending the used range from R202 to R203), then
# populate the newly inserted row 135 with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(135).Insert()

$ws.Cells.Item(135, 1).Value = 4
$ws.Cells.Item(135, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(135, 3).Value = "Los Lagos"
$ws.Cells.Item(135, 4).Value = 44572
$ws.Cells.Item(135, 5).Value = 10
$ws.Cells.Item(135, 6).Value = 100112003
$ws.Cells.Item(135, 7).Value = "Ajo"
$ws.Cells.Item(135, 8).Value = "Chino"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 240
$ws.Cells.Item(135, 11).Value = 20000
$ws.Cells.Item(135, 12).Value = 22000
$ws.Cells.Item(135, 13).Value = 21000
$ws.Cells.Item(135, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(135, 15).Value = "China"
$ws.Cells.Item(135, 16).Value = 2100
$ws.Cells.Item(135, 17).Value = 10
$ws.Cells.Item(135, 18).Value = "Hortaliza"
